$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.240.20"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.849.42"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.50"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6991"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07732"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3063"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.58"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.44"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "1.845.78"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.135"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6869"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.637"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008324"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "29.201.99"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.53"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "2.084.27"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.77"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.527"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1520"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.27"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.836"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.541"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.239"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.191"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.197"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7924"
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.876"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.152"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "1.316.54"
$ws.Range("E38").Value = "  +7.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01875"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9504"
$ws.Range("E41").Value = "  +5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.053"
$ws.Range("E42").Value = "  +6.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.70"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.759"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "1.986.26"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5181"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.25"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.768"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.004"
$ws.Range("E51").Value = "  -0.41%  "
